$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title
$d.Content.Find.Execute("(Name of your distributed application)", $true, $false, $false, $false, $false, $true, 1, $false, "Word Scramble", 2) | Out-Null

# 2. Subtitle
$d.Content.Find.Execute("Architectural Design", $true, $false, $false, $false, $false, $true, 1, $false, "Overview", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Paragraph 4: "This should be an updated version..." (first occurrence,
#    right after "Introduction and Context") turns into the Word Scramble
#    overview paragraph, loses its italics, and gains sz24 + double spacing.
$p4 = $d.Paragraphs(4)
$p4.Range.Text = "Our distributed application will be a two-player game called Word Scramble.  It will largely be based off the board game Scrabble but will have some key differences in game play.  Each player will be given an allotted number of letters and will take turns spelling out words on a 16 x 16 tile game board.  When a player places a word on the board, that word is checked against a dictionary of valid words.  If the word is valid, then the turn is passed to the other player.  There is no restriction on where a new word can be played on the game board.  Play continues until the game board is filled and no new words can be played.  The player with the least number of remaining letters is the winner.  "
$p4.Range.Font.Italic = 0
$p4.Range.Font.Size = 12
$p4.Range.Font.SizeBi = 12
$p4.Range.ParagraphFormat.LineSpacingRule = 2
$p4.Range.ParagraphFormat.LeftIndent = 18
$p4.Range.ParagraphFormat.FirstLineIndent = 18

# ---------------------------------------------------------------------------
# 4. New paragraph inserted right after paragraph 4, carrying the second half
#    of the overview (same formatting: sz24, double spacing, indents).
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs(5)
$p5.Range.Text = " This application will allow us to meet the requirements for this project in the following ways:  The application will consist of three distinct shared resources, utilizing a database for the dictionary and to store user information, a file to keep record of games played, and a shared object in the form of a game board.  We will have a Database Server, a Referee Server, and a client program that will communicate and allow for seamless game play.  There will be several communication protocols in the application using both TCP/IP and UDP protocols.  The scope of this application will challenge our skills and abilities but will be attainable within the time frame we have"
$p5.Range.Font.Italic = 0
$p5.Range.Font.Size = 12
$p5.Range.Font.SizeBi = 12
$p5.Range.ParagraphFormat.LineSpacingRule = 2
$p5.Range.ParagraphFormat.LeftIndent = 18
$p5.Range.ParagraphFormat.FirstLineIndent = 18

# ---------------------------------------------------------------------------
# 5. Heading "Users and their Goals" -> "Actors and their Goals" (now
#    paragraph 6 after the insertion above). The two runs collapse into one.
$p6 = $d.Paragraphs(6)
$rng6 = $p6.Range
$rng6.MoveEnd(1, -1) | Out-Null
$rng6.Text = "Actors and their Goals"

# ---------------------------------------------------------------------------
# 6. Paragraph 7 (the two-run italic paragraph right after the heading above)
#    becomes the detailed actors description, made up of two runs (sz24).
$p7 = $d.Paragraphs(7)
$rngText7 = $p7.Range
$rngText7.MoveEnd(1, -1) | Out-Null
$rngText7.Text = "We have identified three actors that have separate goals. The Actors are Users, Referee server, and Database Server. The Users goal will be able to create an account. After account creation the User will be able to Log In to the client. After the User has logged in the User will request a new game. Inside the game a User will be able to place letters, request new letters, submit a word, send a heartbeat, reply to a heartbeat, and end the game. The Database server will serve as the dictionary database, the user database, played game database, and contain a register of Referee Servers. The Database Server will be able to create a dictionary, send a heartbeat, reply to heartbeats, receive information from a referee, update user information, receive words, search for words, and reply to say if the word is allowed, and update played game results. The final actor is the Referee server. The Referee Server will be sending heartbeats to both the Database Server, as well as the User Client. If the User does not respond it ends the game and notifies the second player. If the Database server is not replying it temporarily pauses all games, it is overseeing and tries to reconnect. After a certain time, it will throw an error and end the game.  When the Referee receives a new game request, it will attempt to connect two players together. After it connects to both players it will start to send game update messages.  After a game exits it will send an update to the Database Server. Finally, when the server is launched it will send its information to the Database Server so the Database Server can sed it to a User Client."

# Turn off italic across the whole paragraph (incl. mark) to match the
# target's non-italic paragraph mark, then size the run text only (sz24).
$p7.Range.Font.Italic = 0
$rngSize7 = $p7.Range
$rngSize7.MoveEnd(1, -1) | Out-Null
$rngSize7.Font.Size = 12
$rngSize7.Font.SizeBi = 12

$p7.Range.ParagraphFormat.LineSpacingRule = 2
$p7.Range.ParagraphFormat.LeftIndent = 18
# firstLine indent of 360 twips = 18pt already present from original ind firstLine=360

foreach ($p in $d.Paragraphs) {
    Write-Host "---"
    Write-Host $p.Range.Text
}
